# Auto-generated edit script applying numeric updates to the Lamia_Profits analysis sheets
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) per the scheduled runner's repricing pass.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 4383.121
$ws.Range("I92").Value = 3342.1853
$ws.Range("J92").Value = 9067.333000000001
$ws.Range("K92").Value = 3342.1853
$ws.Range("L92").Value = 9067.333000000001
$ws.Range("M92").Value = -2094.1853
$ws.Range("N92").Value = -11563.333

$ws.Range("H111").Value = 3939.7144
$ws.Range("I111").Value = 2315.8
$ws.Range("K111").Value = 6947.400000000001
$ws.Range("M111").Value = -3880.400000000001

$ws.Range("H134").Value = 69996.664
$ws.Range("J134").Value = 69996.664
$ws.Range("L134").Value = 69996.664
$ws.Range("N134").Value = -80136.664

$ws.Range("H138").Value = 4869.0586
$ws.Range("I138").Value = 5499.4
$ws.Range("J138").Value = 4760.3794
$ws.Range("K138").Value = 16498.2
$ws.Range("L138").Value = 14281.1382
$ws.Range("M138").Value = -11358.2
$ws.Range("N138").Value = -24561.1382

$ws.Range("H141").Value = 2768.3635
$ws.Range("I141").Value = 2248.5293
$ws.Range("J141").Value = 4535.8
$ws.Range("K141").Value = 6745.5879
$ws.Range("L141").Value = 13607.4
$ws.Range("M141").Value = -1565.5879
$ws.Range("N141").Value = -23967.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 4327
$ws.Range("J12").Value = 4327
$ws.Range("L12").Value = 4327
$ws.Range("N12").Value = -4673

$ws.Range("H61").Value = 8390.6
$ws.Range("I61").Value = 7418.5
$ws.Range("J61").Value = 22000
$ws.Range("K61").Value = 7418.5
$ws.Range("L61").Value = 22000
$ws.Range("M61").Value = -7206.5
$ws.Range("N61").Value = -22424

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null

$ws.Range("H97").Value = 905.125
$ws.Range("I97").Value = 824.9048
$ws.Range("K97").Value = 824.9048
$ws.Range("M97").Value = -328.9048

$ws.Range("H136").Value = 8390.6
$ws.Range("I136").Value = 7418.5
$ws.Range("J136").Value = 22000
$ws.Range("K136").Value = 22255.5
$ws.Range("L136").Value = 66000
$ws.Range("M136").Value = -19705.5
$ws.Range("N136").Value = -71100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1299.5454
$ws.Range("I94").Value = 1141.0526
$ws.Range("K94").Value = 1141.0526
$ws.Range("M94").Value = -690.0526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24822.826
$ws.Range("I31").Value = 2669.1667
$ws.Range("J31").Value = 43811.68
$ws.Range("K31").Value = 2669.1667
$ws.Range("L31").Value = 43811.68
$ws.Range("M31").Value = -2374.1667
$ws.Range("N31").Value = -44401.68

$ws.Range("H34").Value = 24822.826
$ws.Range("I34").Value = 2669.1667
$ws.Range("J34").Value = 43811.68
$ws.Range("K34").Value = 2669.1667
$ws.Range("L34").Value = 43811.68
$ws.Range("M34").Value = -2467.1667
$ws.Range("N34").Value = -44215.68

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = $null
$ws.Range("N45").Value = $null

$ws.Range("H58").Value = 4224.325
$ws.Range("I58").Value = 2236.6538
$ws.Range("J58").Value = 7915.7144
$ws.Range("K58").Value = 2236.6538
$ws.Range("L58").Value = 7915.7144
$ws.Range("M58").Value = -2033.6538
$ws.Range("N58").Value = -8321.714400000001

$ws.Range("H59").Value = 30000
$ws.Range("I59").Value = 30000
$ws.Range("K59").Value = 30000
$ws.Range("M59").Value = -28855

$ws.Range("H97").Value = 50000
$ws.Range("J97").Value = 50000
$ws.Range("L97").Value = 50000
$ws.Range("N97").Value = -51982

$ws.Range("H109").Value = 62583.668
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 62583.668
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 62583.668
$ws.Range("M109").Value = $null
$ws.Range("N109").Value = -64663.668

$ws.Range("H122").Value = 5870.615
$ws.Range("I122").Value = 2913.7896
$ws.Range("J122").Value = 13896.286
$ws.Range("K122").Value = 8741.3688
$ws.Range("L122").Value = 41688.858
$ws.Range("M122").Value = -6291.3688
$ws.Range("N122").Value = -46588.858

$ws.Range("H136").Value = 4224.325
$ws.Range("I136").Value = 2236.6538
$ws.Range("J136").Value = 7915.7144
$ws.Range("K136").Value = 6709.9614
$ws.Range("L136").Value = 23747.1432
$ws.Range("M136").Value = -4159.9614
$ws.Range("N136").Value = -28847.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5967.407
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 6588.3335
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 19765.0005
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -21387.0005

$ws.Range("H71").Value = 5967.407
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 6588.3335
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 59295.0015
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -67407.0015

$ws.Range("H88").Value = 18004
$ws.Range("J88").Value = 19005.334
$ws.Range("L88").Value = 57016.00199999999
$ws.Range("N88").Value = -57872.00199999999

$ws.Range("H91").Value = 18004
$ws.Range("J91").Value = 19005.334
$ws.Range("L91").Value = 57016.00199999999
$ws.Range("N91").Value = -59980.00199999999

$ws.Range("H117").Value = 71433600
$ws.Range("I117").Value = 2376.3333
$ws.Range("J117").Value = 125007016
$ws.Range("K117").Value = 7128.999899999999
$ws.Range("L117").Value = 375021048
$ws.Range("M117").Value = -3686.999899999999
$ws.Range("N117").Value = -375027932

$ws.Range("H129").Value = 5557590.5
$ws.Range("I129").Value = 644.4
$ws.Range("K129").Value = 1933.2
$ws.Range("M129").Value = 3066.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 1456
$ws.Range("I4").Value = 1456
$ws.Range("K4").Value = 1456
$ws.Range("M4").Value = -1344

$ws.Range("H69").Value = 35000
$ws.Range("J69").Value = 35000
$ws.Range("L69").Value = 35000
$ws.Range("N69").Value = -36498

$ws.Range("H72").Value = 35000
$ws.Range("J72").Value = 35000
$ws.Range("L72").Value = 105000
$ws.Range("N72").Value = -112488

$ws.Range("H80").Value = 5550.273
$ws.Range("I80").Value = 3502.7856
$ws.Range("K80").Value = 3502.7856
$ws.Range("M80").Value = -2504.7856

$ws.Range("H83").Value = 5550.273
$ws.Range("I83").Value = 3502.7856
$ws.Range("K83").Value = 17513.928
$ws.Range("M83").Value = -12521.928

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null

$ws.Range("H132").Value = 4513.5884
$ws.Range("I132").Value = 3981.5
$ws.Range("J132").Value = 6996.6665
$ws.Range("K132").Value = 11944.5
$ws.Range("L132").Value = 20989.9995
$ws.Range("M132").Value = -9414.5
$ws.Range("N132").Value = -26049.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 16201.379
$ws.Range("I93").Value = 9273.913
$ws.Range("J93").Value = 42756.668
$ws.Range("K93").Value = 9273.913
$ws.Range("L93").Value = 42756.668
$ws.Range("M93").Value = -8025.913
$ws.Range("N93").Value = -45252.668

$ws.Range("H117").Value = 73684
$ws.Range("J117").Value = 73684
$ws.Range("L117").Value = 73684
$ws.Range("N117").Value = -82862

$ws.Range("H122").Value = 9639
$ws.Range("I122").Value = 5553
$ws.Range("K122").Value = 16659
$ws.Range("M122").Value = -14209

$ws.Range("H132").Value = 5383.4736
$ws.Range("I132").Value = 3948.077
$ws.Range("K132").Value = 11844.231
$ws.Range("M132").Value = -9314.231

$ws.Range("H136").Value = 7763.241
$ws.Range("J136").Value = 10245.533
$ws.Range("L136").Value = 30736.599
$ws.Range("N136").Value = -35836.599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16262

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = $null
$ws.Range("N52").Value = $null

$ws.Range("H57").Value = 54997.5
$ws.Range("I57").Value = 54997.5
$ws.Range("K57").Value = 54997.5
$ws.Range("M57").Value = -54243.5

$ws.Range("H136").Value = 2716.1714
$ws.Range("I136").Value = 1521.6538
$ws.Range("K136").Value = 4564.9614
$ws.Range("M136").Value = -2014.9614
